# Apply the "st 21. 04. 2021" daily-stats update: refreshed AgTests (F) /
# AgPosit (G) figures for the affected rows of the Slovakia Covid daily
# stats sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    362 = @{ F = 229495; G = 3187 }
    363 = @{ F = 188692 }
    364 = @{ F = 168400; G = 2484 }
    365 = @{ F = 184976; G = 2398 }
    366 = @{ F = 339395 }
    367 = @{ F = 766231; G = 3919 }
    368 = @{ F = 346891 }
    369 = @{ F = 234741; G = 2604 }
    370 = @{ F = 180713; G = 2045 }
    371 = @{ F = 160126; G = 1960 }
    372 = @{ F = 179198; G = 1855 }
    373 = @{ F = 349864; G = 2375 }
    374 = @{ F = 773400; G = 3421 }
    375 = @{ F = 351849; G = 1859 }
    376 = @{ F = 221287; G = 2224 }
    377 = @{ F = 176991 }
    378 = @{ F = 157258; G = 1550 }
    381 = @{ F = 746161; G = 2693 }
    382 = @{ F = 356780; G = 1565 }
    385 = @{ F = 150902 }
    386 = @{ F = 182502 }
    391 = @{ F = 176188 }
    392 = @{ F = 220981 }
    393 = @{ F = 307289 }
    395 = @{ F = 749810; G = 1954 }
    396 = @{ F = 164476 }
    398 = @{ F = 298211; G = 1469 }
    399 = @{ F = 200608 }
    400 = @{ F = 150122; G = 757 }
    402 = @{ F = 715546; G = 1380 }
    403 = @{ F = 350402; G = 730 }
    406 = @{ F = 170320; G = 678 }
    407 = @{ F = 156870; G = 670 }
    408 = @{ F = 300974; G = 834 }
    409 = @{ F = 681493; G = 974 }
    410 = @{ F = 345455; G = 610 }
    411 = @{ F = 222468; G = 818 }
}

foreach ($row in ($updates.Keys | Sort-Object)) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
